# This script updates the "dSF" column (column F) values for a set of rows
# in the active worksheet, per the data repull / mean recalculation described
# in the commit message. Only the F-column cells listed below change value;
# all other cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = 2
    10 = 4
    11 = 0
    14 = 1
    18 = 3
    25 = -1
    26 = 0
    30 = -2
    38 = -1
    41 = -4
    49 = -2
    51 = -2
    52 = -2
    53 = 1
    55 = 2
    58 = -2
    60 = 1
    65 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
